$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.297.99'
$ws.Range('E2').Value = '  -1.22%  '

$ws.Range('D3').Value = '3.692.18'
$ws.Range('E3').Value = '  -2.77%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.13%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.46'
$ws.Range('E5').Value = '  +0.16%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.96'
$ws.Range('E6').Value = '  -3.41%  '

$ws.Range('D7').Value = '3.693.65'
$ws.Range('E7').Value = '  -2.56%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.534'
$ws.Range('E9').Value = '  -0.09%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.164'
$ws.Range('E10').Value = '  +2.94%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.19'
$ws.Range('E11').Value = '  -2.07%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.459'
$ws.Range('E12').Value = '  -2.32%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.78'
$ws.Range('E13').Value = '  -1.59%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000242'
$ws.Range('E14').Value = '  -0.81%  '

$ws.Range('D15').Value = '4.305.41'
$ws.Range('E15').Value = '  -2.84%  '

$ws.Range('D16').Value = '3.690.73'
$ws.Range('E16').Value = '  -2.81%  '

$ws.Range('D17').Value = '67.265.15'
$ws.Range('E17').Value = '  -1.49%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.25'
$ws.Range('E18').Value = '  -0.23%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.115'
$ws.Range('E19').Value = '  -1.09%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.23'
$ws.Range('E20').Value = '  +6.83%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '485.35'
$ws.Range('E21').Value = '  -0.73%  '

$ws.Range('E22').Value = '  -2.14%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.721'
$ws.Range('E23').Value = '  -2.16%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.66'
$ws.Range('E24').Value = '  -1.66%  '

$ws.Range('E25').Value = '  +1.57%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.29'
$ws.Range('E26').Value = '  -4.04%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.21'
$ws.Range('E27').Value = '  -0.56%  '

$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.04'
$ws.Range('E28').Value = '  -2.09%  '

$ws.Range('B29').Value = 'Dai'
$ws.Range('C29').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.12%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.91'
$ws.Range('E30').Value = '  -1.14%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.35'
$ws.Range('E31').Value = '  -3.90%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.65'
$ws.Range('E32').Value = '  +0.28%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.21'
$ws.Range('E33').Value = '  -4.01%  '

$ws.Range('D34').Value = '3.827.16'
$ws.Range('E34').Value = '  -2.79%  '

$ws.Range('E35').Value = '  -2.43%  '

$ws.Range('D36').Value = '3.632.38'
$ws.Range('E36').Value = '  -2.77%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.998'
$ws.Range('E37').Value = '  -0.10%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.997'
$ws.Range('E38').Value = '  -2.18%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.82'
$ws.Range('E39').Value = '  -0.53%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.132'
$ws.Range('E40').Value = '  -2.91%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.321'
$ws.Range('E41').Value = '  -1.51%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '48.78'
$ws.Range('E42').Value = '  -0.74%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '426.09'
$ws.Range('E43').Value = '  -5.52%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.94'
$ws.Range('E44').Value = '  -3.77%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.80'
$ws.Range('E45').Value = '  -2.41%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.40'
$ws.Range('E46').Value = '  +0.59%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.31'
$ws.Range('E48').Value = '  -2.82%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.31'
$ws.Range('E49').Value = '  +2.75%  '

$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0350'
$ws.Range('E50').Value = '  -1.02%  '

$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Value = '2.742.01'
$ws.Range('E51').Value = '  -4.09%  '
